$d = $word.ActiveDocument

# Locate the two adjacent empty paragraphs that follow the "...ZBrush(R)."
# paragraph (they share the same paragraph formatting: rFonts
# asciiTheme/hAnsiTheme="minorBidi" + lang="en-US"). The new paragraph of
# text belongs between them.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r" -and $i -lt $d.Paragraphs.Count) {
        $prev = $d.Paragraphs.Item($i - 1)
        if ($prev.Range.Text -like "*ZBrush*") {
            $target = $i
            break
        }
    }
}

# $target now indexes the first of the two empty paragraphs
# (paraId 0987FD02 in the original markup); the new paragraph is inserted
# right after it, before the second empty paragraph (paraId 483CC69C).
$secondEmpty = $d.Paragraphs.Item($target + 1)
$rng = $secondEmpty.Range
$rng.Collapse(1)   # wdCollapseStart

$newText = "Shaders and materials define the appearance of your objects when light interacts with them and are critical to expressing the style of your environment, whether it is intended to realistic or artistic.`r"
$rng.InsertBefore($newText)
